$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item(1).Name = "growth-restored cell lineage"
$wb.Worksheets.Item(2).Name = "growth-halted cell lineage"
$wb.Worksheets.Item(3).Name = "non-deleted cell lineage"
